$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
# (shifts old N/O/P "Late"/heading/"Outstanding" columns one to the right)
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 9.83

# --- Transactions sheet was previously the selected tab; move off of it ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()
$wsTrans.Range("C2").Select() | Out-Null

# Repayment schedule becomes the workbook's active tab, selection at R6
$ws.Activate()
$ws.Range("R6").Select() | Out-Null
